$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 2-22 (A, B, C, D, E)
$data = @(
    @{Row=2;  A=6;   B=3; C=50; D="3,5,6";  E="30,19,1"},
    @{Row=3;  A=20;  B=3; C=50; D="6,15";   E="19,23"},
    @{Row=4;  A=74;  B=3; C=50; D="4,8,14"; E="2,20,26"},
    @{Row=5;  A=120; B=1; C=30; D="10";     E="29"},
    @{Row=6;  A=154; B=3; C=50; D="20";     E="14"},
    @{Row=7;  A=162; B=3; C=50; D="1,7";    E="25,25"},
    @{Row=8;  A=180; B=3; C=50; D="13";     E="33"},
    @{Row=9;  A=266; B=1; C=30; D="21";     E="28"},
    @{Row=10; A=329; B=3; C=50; D="12,18";  E="46,4"},
    @{Row=11; A=381; B=3; C=50; D="11";     E="16"},
    @{Row=12; A=409; B=2; C=40; D="19";     E="39"},
    @{Row=13; A=457; B=2; C=40; D="9,11";   E="30,10"},
    @{Row=14; A=498; B=2; C=40; D="26";     E="8"},
    @{Row=15; A=594; B=3; C=50; D="18,23";  E="12,18"},
    @{Row=16; A=607; B=3; C=50; D="17,22";  E="36,14"},
    @{Row=17; A=633; B=3; C=50; D="27";     E="50"},
    @{Row=18; A=647; B=3; C=50; D="16";     E="27"},
    @{Row=19; A=785; B=1; C=30; D="22,24";  E="6,5"},
    @{Row=20; A=843; B=2; C=40; D="25";     E="40"},
    @{Row=21; A=908; B=1; C=30; D="24";     E="30"},
    @{Row=22; A=955; B=3; C=50; D="2,4";    E="24,26"}
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.A
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C

    # Prefix with an apostrophe to force text entry for values that would
    # otherwise look like a plain number (e.g. "10"), then restore the
    # cell's original (Normal/default) style so no stray number format is
    # left behind on the cell.
    $dCell = $ws.Cells.Item($r, 4)
    $dCell.Value = "'" + $item.D
    $dCell.Style = "Normal"

    $eCell = $ws.Cells.Item($r, 5)
    $eCell.Value = "'" + $item.E
    $eCell.Style = "Normal"
}
